$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12  = 4.951000000000001
    32  = 6.77
    36  = 8.599
    38  = 6.149999999999999
    46  = 6.248
    54  = 4.933
    55  = 4.763
    67  = 5.505
    69  = 5.367000000000001
    72  = 5.7
    91  = 5.564
    99  = 5.245
    104 = 7.950999999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
